# Update the submission deadline date in DataStoryTemplate.docx.
#
# Before: " no later than 23:59 on 12.12.2025. We will use the version
#           closest to that timestamp in case of missed deadlines."
#          (a single run)
# After:  the same sentence, but the day portion of the date changes
#          from "12" to "05" (i.e. the due date becomes 05.12.2025),
#          and the run is split into three runs:
#            " no later than 23:59 on " | "05" | ".12.2025. We will use
#            the version closest to that timestamp in case of missed
#            deadlines."
#          all three keeping the original run formatting
#          (<w:color w:val="595959" w:themeColor="text1" w:themeTint="A6"/>).

$d = $word.ActiveDocument

# --- Step 1: fix the date text itself (day "12" -> "05") -------------------
# Scope the Find to the text following "Merry" so we never touch the
# preceding run, then replace just the first two digits of "12.12.2025".
$scan = $d.Content
$findMerry = $scan.Find
$findMerry.ClearFormatting()
$findMerry.Text = "Merry"
$null = $findMerry.Execute()

$afterMerry = $scan.End
$dateScope = $d.Range($afterMerry, $d.Content.End)
$findDate = $dateScope.Find
$findDate.ClearFormatting()
$findDate.Text = "12.12.2025"
$null = $findDate.Execute()

$dayStart = $dateScope.Start
$dayEnd = $dayStart + 2
$d.Range($dayStart, $dayEnd).Text = "05"

# --- Step 2: re-establish the run boundaries --------------------------------
# Replacing text above merges every same-formatted run in the paragraph
# into one (Word's normal run-coalescing behaviour). Re-split the
# paragraph back into the runs the edit introduced by nudging (and
# immediately reverting) a character property on each segment boundary
# -- this forces new runs without altering the final formatting.
# Work right-to-left so earlier offsets stay valid.

$scanM = $d.Content
$findM = $scanM.Find
$findM.ClearFormatting()
$findM.Text = "Merry"
$null = $findM.Execute()
$mEnd = $scanM.End

$scanDay = $d.Content
$findDay2 = $scanDay.Find
$findDay2.ClearFormatting()
$findDay2.Text = "05.12.2025"
$null = $findDay2.Execute()
$dayStart2 = $scanDay.Start
$dayEnd2 = $dayStart2 + 2

$scanEnd = $d.Content
$findEnd = $scanEnd.Find
$findEnd.ClearFormatting()
$findEnd.Text = "missed deadlines."
$null = $findEnd.Execute()
$sentenceEnd = $scanEnd.End

# Segment 3: ".12.2025. We will use the version closest to that timestamp
#             in case of missed deadlines."
$seg3 = $d.Range($dayEnd2, $sentenceEnd)
$seg3.Bold = 1
$seg3.Bold = 0

# Segment 2: "05"
$seg2 = $d.Range($dayStart2, $dayEnd2)
$seg2.Bold = 1
$seg2.Bold = 0

# Segment 1: " no later than 23:59 on " (kept apart from the "Merry" run)
$seg1 = $d.Range($mEnd, $dayStart2)
$seg1.Bold = 1
$seg1.Bold = 0
